# "Add damping. Fix shoulder joints"
#
# Joint Limits sheet:
#  - Hip pitch Lower/Upper values were swapped (damping/sign fix) for
#    left_hip (rows 6-7) and right_hip (rows 8-9).
#  - Knee (rows 10-11) and the second knee block (rows 12-13) pitch
#    Lower/Upper values swap too, and a clarifying note is added to
#    row 12 ("Hip cannot go below 0").
#  - Shoulder joints (rows 20-23) had wrong pitch/yaw magnitudes and a
#    yaw formula that referenced the wrong row; both are corrected, and
#    the "Ignoring rotation" note on the left shoulder is replaced with
#    a TODO about adding rotation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Joint Limits")
$ws.Activate()

# --- left_hip / right_hip: Pitch (Degrees) Lower<->Upper swap ---
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = 123.8
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 123.8

# --- left_knee: Pitch (Degrees) Lower<->Upper swap ---
$ws.Range("C10").Value = -117.1
$ws.Range("C11").Value = 0

# --- right_knee: Pitch (Degrees) Lower<->Upper swap + new note ---
$ws.Range("C12").Value = -117.1
$ws.Range("C13").Value = 0
$ws.Range("I12").Value = "Hip cannot go below 0"

# --- left_shoulder / right_shoulder: fix Pitch/Yaw magnitudes and the
#     Yaw (Rad) formula (it was pointing at the wrong row) ---
$ws.Range("C20").Value = -193.2
$ws.Range("D20").Value = -132.1
$ws.Range("F20").Formula = "=RADIANS(D20)"

$ws.Range("C21").Value = 63
$ws.Range("D21").Value = 50.8
$ws.Range("F21").Formula = "=RADIANS(D21)"

$ws.Range("C22").Value = -193.2
$ws.Range("D22").Value = -132.1
$ws.Range("F22").Formula = "=RADIANS(D22)"

$ws.Range("C23").Value = 63
$ws.Range("D23").Value = 50.8
$ws.Range("F23").Formula = "=RADIANS(D23)"

$ws.Range("I20").Value = "TODO: Add rotation"

# --- view state ---
$ws.Range("I21").Select()

$win = $excel.ActiveWindow
$win.Left = 2700
$win.Top = 2600
